# Auto-generated edit script applying the crypto price update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

$ws.Range("D2").Value = "26.297.23"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.601.09"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.12%  "
Set-TextCell "D5" "212.55"
$ws.Range("E5").Value = "  +0.33%  "
Set-TextCell "D6" "0.502"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +0.04%  "
Set-TextCell "D10" "18.99"
$ws.Range("E10").Value = "  -1.38%  "
Set-TextCell "D11" "0.0856"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "1.826.82"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "1.597.15"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  -2.10%  "
Set-TextCell "D16" "63.70"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "26.293.80"
$ws.Range("E17").Value = "  +0.29%  "
Set-TextCell "D18" "230.07"
$ws.Range("E18").Value = "  +6.92%  "
$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("E19").Value = "  -0.56%  "
Set-TextCell "D20" "7.61"
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("E21").Value = "  +0.08%  "
Set-TextCell "D22" "4.28"
$ws.Range("E22").Value = "  +0.99%  "
Set-TextCell "D23" "2.16"
$ws.Range("E23").Value = "  -0.71%  "
Set-TextCell "D24" "8.95"
$ws.Range("E24").Value = "  -0.21%  "
Set-TextCell "D25" "145.63"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +0.17%  "
Set-TextCell "D27" "6.96"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  +1.09%  "
Set-TextCell "D29" "15.43"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  +0.83%  "
Set-TextCell "D32" "3.20"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "1.446.33"
$ws.Range("E33").Value = "  +6.15%  "
Set-TextCell "D34" "2.96"
$ws.Range("E34").Value = "  +0.74%  "
Set-TextCell "D35" "2.42"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  +0.58%  "
Set-TextCell "D37" "0.570"
$ws.Range("E37").Value = "  -1.94%  "
$ws.Range("E38").Value = "  -1.11%  "
Set-TextCell "D39" "0.823"
$ws.Range("E39").Value = "  +0.57%  "
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D43" "0.922"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.739.25"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("E45").Value = "  -1.23%  "
Set-TextCell "D46" "60.74"
$ws.Range("E46").Value = "  -0.40%  "
Set-TextCell "D47" "87.38"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D48" "1.49"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D49" "0.0500"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D50" "0.0950"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextCell "D51" "0.998"
$ws.Range("E51").Value = "  +0.13%  "
